# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.116.85'
$ws.Range("E2").Value = '  -0.39%  '
$ws.Range("D3").Value = '2.454.51'
$ws.Range("E3").Value = '  -3.13%  '
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '525.43'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.21%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '131.06'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -1.69%  '
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.566'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +0.28%  '
$ws.Range("D9").Value = '2.460.71'
$ws.Range("E9").Value = '  -2.79%  '
$ws.Range("E10").Value = '  -0.37%  '
$ws.Range("E11").Value = '  -1.95%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.00'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -3.48%  '
$ws.Range("E13").Value = '  -1.96%  '
$ws.Range("D14").Value = '2.888.11'
$ws.Range("E14").Value = '  -2.95%  '
$ws.Range("D15").Value = '58.034.92'
$ws.Range("E15").Value = '  -0.54%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.72'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -2.30%  '
$ws.Range("E17").Value = '  -1.57%  '
$ws.Range("D18").Value = '2.456.95'
$ws.Range("E18").Value = '  -2.68%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.46'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -2.07%  '
$ws.Range("E20").Value = '  -1.48%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '315.93'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -2.25%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.15'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.12%  '
$ws.Range("E23").Value = '  -0.10%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.43'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.67%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.405'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.69%  '
$ws.Range("E26").Value = '  +0.32%  '
$ws.Range("D27").Value = '2.567.88'
$ws.Range("E27").Value = '  -2.29%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.158'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -1.48%  '
$ws.Range("E29").Value = '  -1.49%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '173.75'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +3.14%  '
$ws.Range("E31").Value = '  -1.61%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.70'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -1.62%  '
$ws.Range("E33").Value = '  -2.09%  '
$ws.Range("E34").Value = '  -4.28%  '
$ws.Range("E35").Value = '  -0.02%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.997'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.11%  '
$ws.Range("E37").Value = '  -2.04%  '
$ws.Range("E38").Value = '  -5.27%  '
$ws.Range("E39").Value = '  -3.26%  '
$ws.Range("B40").Value = 'SuiNetwork'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.813'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +4.95%  '
$ws.Range("B41").Value = 'OKB'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '36.29'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.37%  '
$ws.Range("E42").Value = '  -2.35%  '
$ws.Range("E43").Value = '  -1.35%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '262.62'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -5.77%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.588'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -2.69%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.81'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -4.21%  '
$ws.Range("E47").Value = '  +0.52%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '122.34'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -5.52%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0495'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -1.17%  '
$ws.Range("E50").Value = '  -1.34%  '
$ws.Range("E51").Value = '  -3.94%  '
